$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update price list rows (A2:D16) -----------------------------------
# Row 2: Arduino mega 2560 pro mini
$ws.Range("A2").Value = "Arduino mega 2560 pro mini"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1021
$ws.Range("D2").Value = 468.5

# Row 3: Arduino nano
$ws.Range("A3").Value = "Arduino nano"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 284.05
$ws.Range("D3").Value = 107.5

# Row 4: Wemos d1 mini
$ws.Range("A4").Value = "Wemos d1 mini"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 270.75799999999998
$ws.Range("D4").Value = 119.35

# Row 5: Дисплей Nextion
$ws.Range("A5").Value = "Дисплей Nextion"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 2355
$ws.Range("D5").Value = 800

# Row 6: Модуль SD карты
$ws.Range("A6").Value = "Модуль SD карты "
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 89
$ws.Range("D6").Value = 27.94

# Row 7: SD карта (previously blank)
$ws.Range("A7").Value = "SD карта"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 250
$ws.Range("D7").Value = 83.8

# Row 8: Шаговый двигатель 28byj-48 (previously blank)
$ws.Range("A8").Value = "Шаговый двигатель 28byj-48"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 240
$ws.Range("D8").Value = 92.78

# Row 9: Драйвер двигателя ULN2003 (previously blank)
$ws.Range("A9").Value = "Драйвер двигателя ULN2003"
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 50
$ws.Range("D9").Value = 32.380000000000003

# Row 10: Модуль реального времени DS1307 (previously blank)
$ws.Range("A10").Value = "Модуль реального времени DS1307"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 97
$ws.Range("D10").Value = 34.92

# Row 11: Сервопривод mg90s (previously blank)
$ws.Range("A11").Value = "Сервопривод mg90s"
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 160
$ws.Range("D11").Value = 124.56

# Row 12: Пьезо-пищалка (previously blank)
$ws.Range("A12").Value = "Пьезо-пищалка"
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 50
$ws.Range("D12").Value = 8.8000000000000007

# Row 13: Сенсорная кнопка TTP223 (previously blank)
$ws.Range("A13").Value = "Сенсорная кнопка TTP223"
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 25
$ws.Range("D13").Value = 12

# Row 14: Переключатель (previously blank)
$ws.Range("A14").Value = "Переключатель"
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 15
$ws.Range("D14").Value = 3

# Row 15: Замок (previously blank)
$ws.Range("A15").Value = "Замок"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 92.2
$ws.Range("D15").Value = 51.3
$ws.Range("E15").ClearContents()

# Row 16: Стилус (previously blank)
$ws.Range("A16").Value = "Стилус"
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = 137
$ws.Range("D16").Value = 32.380000000000003
$ws.Range("E16").ClearContents()
$ws.Range("F16").ClearContents()

# Row 17 totals stay as formulas; clear the stray F17 placeholder cell
$ws.Range("F17").ClearContents()

# Row 18: keep a formatted-but-empty B18 cell in place (matches the
# author's edit, which left a formatted but value-less cell there)
$ws.Range("B18").Font.Bold = $false

# --- Column width: widen column A to fit the longer product names ------
$ws.Columns("A:A").AutoFit()

# --- View: move the freeze pane up and select the total cell -----------
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
[void]$ws.Range("C17").Select()
